$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update the Attribute column (C) for rows 7-9 from "NA" to "TERTV8"
$ws.Range("C7").Value = "TERTV8"
$ws.Range("C8").Value = "TERTV8"
$ws.Range("C9").Value = "TERTV8"

# Update the selection to match the saved cursor position in the file
$ws.Range("C10").Select()
